$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - copy date style/format from the row above (A19) then set the new date value
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = (Get-Date -Year 2025 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B20").Value = 18
$ws.Range("C20").Value = 50
$ws.Range("D20").Value = 21
$ws.Range("E20").Value = 46
$ws.Range("F20").Value = "CS Introduction Lecture 16"

# Row 21
$ws.Range("A19").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = (Get-Date -Year 2025 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B21").Value = 19
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = "CS Introduction Lecture 16"

$ws.Range("F21").Select()
